{"js": "// Find the run of text that renders the empty-state placeholder for the\n// \"otherAttachments\" table (\"{d.otherAttachments:ifEM():show(.noData)}\")\n// and append the missing \":elseShow()\" clause so the template no longer\n// falls through to printing \"[object object]\" when the table has rows.\nconst results = context.document.body.search(\n  \"otherAttachments:ifEM():show(.noData)\",\n  { matchCase: true, matchWholeWord: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text 'otherAttachments:ifEM():show(.noData)' not found.\");\n}\n\n// There should be exactly one match (the bare field, as opposed to the\n// otherAttachments[i]/[i+1] rows inside the table above it).\nconst target = results.items[0];\ntarget.insertText(\":elseShow()\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Find the run of text that renders the empty-state placeholder for the\n# \"otherAttachments\" table (\"{d.otherAttachments:ifEM():show(.noData)}\")\n# and append the missing \":elseShow()\" clause so the template no longer\n# falls through to printing \"[object object]\" when the table has rows.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"otherAttachments:ifEM():show(.noData)\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWildcards = $false\n$found = $rng.Find.Execute()\n\nif (-not $found) {\n    throw \"Target text 'otherAttachments:ifEM():show(.noData)' not found.\"\n}\n\n# Collapse the found range to its end point so the insertion lands right\n# after \"...show(.noData)\" and before the closing \"}\" run.\n$rng.Collapse(0)\n$rng.InsertAfter(\":elseShow()\")\n"}
